# Update "想去人数" (want-to-go count) figures across the four sheets of
# 广州-漫展信息.xlsx to the freshly scraped values (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 1236
$ws.Range("F4").Value  = 12818
$ws.Range("F5").Value  = 729
$ws.Range("F10").Value = 1867
$ws.Range("F13").Value = 514
$ws.Range("F16").Value = 346
$ws.Range("F19").Value = 130
$ws.Range("F24").Value = 1291
$ws.Range("F25").Value = 335
$ws.Range("F26").Value = 65

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 157
$ws.Range("F8").Value = 13

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 875
$ws.Range("F3").Value = 4077

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 875
$ws.Range("F6").Value  = 1236
$ws.Range("F7").Value  = 12818
$ws.Range("F9").Value  = 729
$ws.Range("F10").Value = 4077
$ws.Range("F15").Value = 1867
$ws.Range("F18").Value = 514
$ws.Range("F21").Value = 157
$ws.Range("F22").Value = 157
$ws.Range("F25").Value = 13
$ws.Range("F29").Value = 346
$ws.Range("F33").Value = 130
$ws.Range("F41").Value = 1291
$ws.Range("F43").Value = 335
$ws.Range("F44").Value = 65

$wb.Save()
